$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update existing row 39 ---
$ws.Range("E39").Value = 6750
$ws.Range("G39").Value = "13-Nov-2025 00:00:00"
$ws.Range("V39").Value = "RS"
$ws.Range("AA39").Value = 5
$ws.Range("AB39").Value = 0.9
$ws.Range("AE39").Value = 5.9

# --- Update existing row 40 ---
$ws.Range("E40").Value = 8350
$ws.Range("G40").Value = "13-Nov-2025 00:00:00"
$ws.Range("V40").Value = "RS"
$ws.Range("AA40").Value = 5
$ws.Range("AB40").Value = 0.9
$ws.Range("AE40").Value = 5.9

# --- Add new row 41 ---
$ws.Range("A41").Value = 16648
$ws.Range("B41").Value = "TALACHINTALA SURYA NIMROD"
$ws.Range("C41").Value = 9492203981
$ws.Range("D41").Value = 7750
$ws.Range("E41").Value = "'"
$ws.Range("F41").Value = "13-Nov-2025 08:16:51"
$ws.Range("G41").Value = "'"
$ws.Range("H41").Value = "TRANSACTION IS SUCCESSFUL"
$ws.Range("I41").Value = "OK"
$ws.Range("J41").Value = "lVl"
$ws.Range("K41").Value = "seven thousand seven hundred fifty"
$ws.Range("L41").Value = 100000036600
$ws.Range("M41").Value = "SALESIAN EDUCATION SOCIETY"
$ws.Range("N41").Value = 753702
$ws.Range("O41").Value = 1234
$ws.Range("P41").Value = 11000316431279
$ws.Range("Q41").Value = 1763001959
$ws.Range("R41").Value = 108562669809
$ws.Range("S41").Value = "INR"
$ws.Range("T41").Value = "sale"
$ws.Range("U41").Value = "ICICI UPI QR"
$ws.Range("V41").Value = "NRNS"
$ws.Range("W41").Value = "SIBL0000899"
$ws.Range("X41").Value = "MERCHANT"
$ws.Range("Y41").Value = "UPI"
$ws.Range("Z41").Value = "kotakschoolvsp@gmail.com"
$ws.Range("AA41").Value = "'"
$ws.Range("AB41").Value = "'"
$ws.Range("AC41").Value = 0
$ws.Range("AD41").Value = 0
$ws.Range("AE41").Value = "'"
$ws.Range("AF41").Value = "KOTAK SALESIAN PRIMARY SCHOOL"
$ws.Range("AG41").Value = "'"
$ws.Range("AH41").Value = "REGULAR"
$ws.Range("AI41").Value = 18172
$ws.Range("AJ41").Value = "'265254"
$ws.Range("AK41").Value = "'2048"
$ws.Range("AL41").Value = "'"

# --- Add new row 42 ---
$ws.Range("A42").Value = 17187
$ws.Range("B42").Value = "CHINTHA VAISHNAVI"
$ws.Range("C42").Value = 8886428971
$ws.Range("D42").Value = 6750
$ws.Range("E42").Value = "'"
$ws.Range("F42").Value = "13-Nov-2025 17:32:13"
$ws.Range("G42").Value = "'"
$ws.Range("H42").Value = "TRANSACTION IS SUCCESSFUL"
$ws.Range("I42").Value = "OK"
$ws.Range("J42").Value = "PREKGUKG"
$ws.Range("K42").Value = "six thousand seven hundred fifty"
$ws.Range("L42").Value = 100000036600
$ws.Range("M42").Value = "SALESIAN EDUCATION SOCIETY"
$ws.Range("N42").Value = 753702
$ws.Range("O42").Value = 1234
$ws.Range("P42").Value = 11000316535819
$ws.Range("Q42").Value = 1763035319
$ws.Range("R42").Value = 173230421882
$ws.Range("S42").Value = "INR"
$ws.Range("T42").Value = "sale"
$ws.Range("U42").Value = "ICICI UPI QR"
$ws.Range("V42").Value = "NRNS"
$ws.Range("W42").Value = "SIBL0000899"
$ws.Range("X42").Value = "MERCHANT"
$ws.Range("Y42").Value = "UPI"
$ws.Range("Z42").Value = "kotakschoolvsp@gmail.com"
$ws.Range("AA42").Value = "'"
$ws.Range("AB42").Value = "'"
$ws.Range("AC42").Value = 0
$ws.Range("AD42").Value = 0
$ws.Range("AE42").Value = "'"
$ws.Range("AF42").Value = "KOTAK SALESIAN SCHOOL MANAGEMENT ACCOUNT"
$ws.Range("AG42").Value = "'"
$ws.Range("AH42").Value = "REGULAR"
$ws.Range("AI42").Value = 19901
$ws.Range("AJ42").Value = "'264986"
$ws.Range("AK42").Value = "'2047"
$ws.Range("AL42").Value = "UPI INTENT"
